$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '29.302.53'
$ws.Cells.Item(2, 5).Value = '  +0.36%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.874.95'
$ws.Cells.Item(3, 5).Value = '  +0.61%  '

# Row 4
$ws.Cells.Item(4, 4).Value = "'1.000"
$ws.Cells.Item(4, 4).Style = 'Normal'
$ws.Cells.Item(4, 5).Value = '  -0.22%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'0.7124"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -0.84%  '

# Row 6
$ws.Cells.Item(6, 4).Value = "'242.43"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +0.75%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +1.18%  '

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.07732"
$ws.Cells.Item(9, 4).Style = 'Normal'
$ws.Cells.Item(9, 5).Value = '  -0.33%  '

# Row 10
$ws.Cells.Item(10, 4).Value = "'25.10"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +0.55%  '

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.08444"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +2.28%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.870.37'
$ws.Cells.Item(12, 5).Value = '  +0.71%  '

# Row 13
$ws.Cells.Item(13, 4).Value = "'5.215"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +0.01%  '

# Row 14
$ws.Cells.Item(14, 4).Value = "'0.7108"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -0.77%  '

# Row 15
$ws.Cells.Item(15, 4).Value = "'91.34"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +1.18%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '29.310.87'
$ws.Cells.Item(16, 5).Value = '  +0.30%  '

# Row 17
$ws.Cells.Item(17, 4).Value = "'0.000008288"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  +6.39%  '

# Row 18
$ws.Cells.Item(18, 4).Value = "'5.984"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  +2.43%  '

# Row 19
$ws.Cells.Item(19, 4).Value = "'242.72"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -0.29%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '2.127.58'
$ws.Cells.Item(20, 5).Value = '  +0.26%  '

# Row 21
$ws.Cells.Item(21, 5).Value = '  +0.63%  '

# Row 22
$ws.Cells.Item(22, 4).Value = "'0.9997"
$ws.Cells.Item(22, 4).Style = 'Normal'

# Row 23
$ws.Cells.Item(23, 4).Value = "'7.811"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -1.57%  '

# Row 24
$ws.Cells.Item(24, 4).Value = "'1.0000"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  -0.27%  '

# Row 25
$ws.Cells.Item(25, 4).Value = "'0.1619"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +1.53%  '

# Row 26
$ws.Cells.Item(26, 4).Value = "'163.21"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  +0.61%  '

# Row 27
$ws.Cells.Item(27, 4).Value = "'9.014"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +1.14%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +1.88%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +0.84%  '

# Row 30
$ws.Cells.Item(30, 4).Value = "'4.419"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +1.53%  '

# Row 31
$ws.Cells.Item(31, 4).Value = "'4.329"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  +6.01%  '

# Row 32
$ws.Cells.Item(32, 4).Value = "'1.285"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -3.25%  '

# Row 33
$ws.Cells.Item(33, 4).Value = "'0.05256"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +1.32%  '

# Row 34
$ws.Cells.Item(34, 4).Value = "'1.922"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +0.44%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'ARBITRUM'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(35, 4).Value = "'1.173"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -0.19%  '

# Row 36
$ws.Cells.Item(36, 2).Value = 'ImmutableX'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(36, 4).Value = "'0.7466"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +2.46%  '

# Row 37
$ws.Cells.Item(37, 4).Value = "'2.685"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +0.24%  '

# Row 38
$ws.Cells.Item(38, 4).Value = "'0.01859"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  +0.59%  '

# Row 39
$ws.Cells.Item(39, 4).Value = "'2.718"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  +0.70%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '1.160.78'
$ws.Cells.Item(40, 5).Value = '  +0.42%  '

# Row 41
$ws.Cells.Item(41, 4).Value = "'6.361"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +4.08%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(42, 4).Value = "'0.8894"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -1.18%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'Aave'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(43, 4).Value = "'72.99"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +1.06%  '

# Row 44
$ws.Cells.Item(44, 4).Value = "'106.49"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +4.83%  '

# Row 45
$ws.Cells.Item(45, 4).Value = "'0.9995"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -0.16%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '2.024.63'
$ws.Cells.Item(46, 5).Value = '  +0.41%  '

# Row 47
$ws.Cells.Item(47, 4).Value = "'1.807"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +2.48%  '

# Row 48
$ws.Cells.Item(48, 4).Value = "'0.5194"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -1.72%  '

# Row 49
$ws.Cells.Item(49, 4).Value = "'0.00000000121"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +4.09%  '

# Row 50
$ws.Cells.Item(50, 4).Value = "'9.394"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +1.06%  '

# Row 51
$ws.Cells.Item(51, 4).Value = "'0.4298"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +1.46%  '
